# Generate Report for handback
# Adds a new handback row (file 29d53477-b708-4a86-8c9a-afacdcac8bc3) as row 4
# to the "Overview", "zh-cn" and "de-de" worksheets, mirroring the existing
# rows 2/3 pattern for files 5fe82f49... and 5c6247bd...

$wb = $excel.ActiveWorkbook

$newFile   = "29d53477-b708-4a86-8c9a-afacdcac8bc3"
$newFileMd = "$newFile.md"
$xlfHash   = "1f5c4d40285c84bdfd58e4101ead0ca81ed574dd"

$statusInSync = "Handed back: in sync with en-US"
$handoffReason = "Include"

function Set-HyperlinkCell($ws, $addr, $url, $display) {
    $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $display)
    $ws.Range($addr).Font.Underline = $true
    $ws.Range($addr).Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a8ebf3a7b936ff53f7e4ea6b775fa56dad46ad51/e2e/$newFileMd"

Set-HyperlinkCell $wsOverview "A4" $overviewMdUrl $newFileMd
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhMdUrl = $overviewMdUrl
$zhHandoffXlf = "$newFile.$xlfHash.zh-cn.xlf"
$zhHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02506fc64a46036a9e148fea4a0d4c2eba9fd687/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhHandoffXlf"
$zhMdUrl2 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/68197b7e922529a956905ae5f7e41462a7d24173/e2e/$newFileMd"
$zhHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3bfc5b6ac8a5a1e4b8d69bcc0943ed11de003aad/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhHandoffXlf"

Set-HyperlinkCell $wsZh "A4" $zhMdUrl $newFileMd
$wsZh.Range("B4").Value = $statusInSync
Set-HyperlinkCell $wsZh "C4" $zhHandoffUrl $zhHandoffXlf
$wsZh.Range("D4").Value = "2016-01-26 09:29:00"
Set-HyperlinkCell $wsZh "E4" $zhMdUrl2 $newFileMd
Set-HyperlinkCell $wsZh "F4" $zhHandbackUrl $zhHandoffXlf
$wsZh.Range("G4").Value = "2016-01-26 09:29:46"
$wsZh.Range("H4").Value = $handoffReason

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deMdUrl = $overviewMdUrl
$deHandoffXlf = "$newFile.$xlfHash.de-de.xlf"
$deHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9c345db84bb85b21a125b3f2f4052309932a4246/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deHandoffXlf"
$deMdUrl2 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/31eeb7c5ef4454af7384aa03c98d8160cdea9122/e2e/$newFileMd"
$deHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2fdd545aaa88057f6eb597bcad032475bd322473/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deHandoffXlf"

Set-HyperlinkCell $wsDe "A4" $deMdUrl $newFileMd
$wsDe.Range("B4").Value = $statusInSync
Set-HyperlinkCell $wsDe "C4" $deHandoffUrl $deHandoffXlf
$wsDe.Range("D4").Value = "2016-01-26 09:29:13"
Set-HyperlinkCell $wsDe "E4" $deMdUrl2 $newFileMd
Set-HyperlinkCell $wsDe "F4" $deHandbackUrl $deHandoffXlf
$wsDe.Range("G4").Value = "2016-01-26 09:30:12"
$wsDe.Range("H4").Value = $handoffReason

Write-Host "Handback report row added for $newFile"
